# Update Nilton's schedule sheet: swap the "MEC-3B-Elem. Máquinas" slot
# from C3 to F4, and swap the "MEC-3A-Metrologia 2" slots from F10/F16
# to F14/F15 (Funcionando com o ELM-2NA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "-"
$ws.Range("F4").Value = "MEC-3B-Elem. Máquinas"

$ws.Range("F10").Value = "-"
$ws.Range("F14").Value = "[-, 'MEC-3A-Metrologia 2', -, -]"
$ws.Range("F15").Value = "[-, 'MEC-3A-Metrologia 2', -, -]"
$ws.Range("F16").Value = "-"
